# Refresh cryptocurrency price/volume/rank data (GitHub Actions scrape update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "236.08"
$cell.Style = "Normal"

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "21.69"
$cell.Style = "Normal"

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "5.363"
$cell.Style = "Normal"

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "0.05569"
$cell.Style = "Normal"

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "3.366"
$cell.Style = "Normal"

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "6.463"
$cell.Style = "Normal"

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.7999"
$cell.Style = "Normal"

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "1.043"
$cell.Style = "Normal"

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.1399"
$cell.Style = "Normal"

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07296"
$cell.Style = "Normal"

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.03188"
$cell.Style = "Normal"

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.02942"
$cell.Style = "Normal"

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.09235"
$cell.Style = "Normal"

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.001666"
$cell.Style = "Normal"

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "3.257"
$cell.Style = "Normal"

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.04759"
$cell.Style = "Normal"

$cell = $ws.Range("B18")
$cell.NumberFormat = "@"
$cell.Value = "One"
$cell.Style = "Normal"

$cell = $ws.Range("C18")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$cell.Style = "Normal"

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.0005712"
$cell.Style = "Normal"

$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "17OneONE"
$cell.Style = "Normal"

$cell = $ws.Range("B19")
$cell.NumberFormat = "@"
$cell.Value = "TigerCash"
$cell.Style = "Normal"

$cell = $ws.Range("C19")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$cell.Style = "Normal"

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.006260"
$cell.Style = "Normal"

$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "18TigerCashTCH"
$cell.Style = "Normal"

$cell = $ws.Range("B20")
$cell.NumberFormat = "@"
$cell.Value = "HotbitToken"
$cell.Style = "Normal"

$cell = $ws.Range("C20")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$cell.Style = "Normal"

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.005066"
$cell.Style = "Normal"

$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "19HotbitTokenHTB"
$cell.Style = "Normal"

$cell = $ws.Range("B21")
$cell.NumberFormat = "@"
$cell.Value = "BitKan"
$cell.Style = "Normal"

$cell = $ws.Range("C21")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$cell.Style = "Normal"

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "0.001048"
$cell.Style = "Normal"

$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "20BitKanKAN"
$cell.Style = "Normal"

$cell = $ws.Range("B22")
$cell.NumberFormat = "@"
$cell.Value = "NitroEx"
$cell.Style = "Normal"

$cell = $ws.Range("C22")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$cell.Style = "Normal"

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.0001499"
$cell.Style = "Normal"

$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "21NitroExNTX"
$cell.Style = "Normal"

$cell = $ws.Range("B23")
$cell.NumberFormat = "@"
$cell.Value = "UpBots"
$cell.Style = "Normal"

$cell = $ws.Range("C23")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$cell.Style = "Normal"

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.0004183"
$cell.Style = "Normal"

$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "22UpBotsUBXT"
$cell.Style = "Normal"

$cell = $ws.Range("B24")
$cell.NumberFormat = "@"
$cell.Value = "LEO"
$cell.Style = "Normal"

$cell = $ws.Range("C24")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$cell.Style = "Normal"

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "3.948"
$cell.Style = "Normal"

$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "23LEOLEOBestin24h"
$cell.Style = "Normal"

$cell = $ws.Range("B25")
$cell.NumberFormat = "@"
$cell.Value = "BTSEToken"
$cell.Style = "Normal"

$cell = $ws.Range("C25")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$cell.Style = "Normal"

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.200"
$cell.Style = "Normal"

$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "24BTSETokenBTSE"
$cell.Style = "Normal"

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.1295"
$cell.Style = "Normal"

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.04120"
$cell.Style = "Normal"

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.007013"
$cell.Style = "Normal"

$cell = $ws.Range("B42")
$cell.NumberFormat = "@"
$cell.Value = "BKEXToken"
$cell.Style = "Normal"

$cell = $ws.Range("C42")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$cell.Style = "Normal"

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.1038"
$cell.Style = "Normal"

$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "41BKEXTokenBKK"
$cell.Style = "Normal"

$cell = $ws.Range("B43")
$cell.NumberFormat = "@"
$cell.Value = "CEJI"
$cell.Style = "Normal"

$cell = $ws.Range("C43")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$cell.Style = "Normal"

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.002919"
$cell.Style = "Normal"

$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "42CEJICEJI"
$cell.Style = "Normal"

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.009505"
$cell.Style = "Normal"

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.00005433"
$cell.Style = "Normal"

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.6803"
$cell.Style = "Normal"

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.03262"
$cell.Style = "Normal"

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.00002101"
$cell.Style = "Normal"
